# Add the next four "Logs" entries documenting work on drag and drop in
# inventory mode (rows 51-54 on the "Logs" sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$entries = @(
    @{ Row = 51; Date = 45559; Text = "add drag and drop first steps to belt, need to add this for inventory/hands and remove object in old place" },
    @{ Row = 52; Date = 45562; Text = "work on drag and drop" },
    @{ Row = 53; Date = 45563; Text = "work on drag and drop ( I hate this)" },
    @{ Row = 54; Date = 45564; Text = "finally, drag and drop works as I want. Items can be merged or change one anothers place. Player cant drop item or split" }
)

foreach ($entry in $entries) {
    $row = $entry.Row
    $prevRow = $row - 1

    # Copy the formatting of the row above down onto the new row so the
    # date column keeps its date number format / alignment and the text
    # column keeps its wrap-text style.
    $srcRange = $ws.Range("A" + $prevRow + ":B" + $prevRow)
    $dstRange = $ws.Range("A" + $row + ":B" + $row)
    $srcRange.Copy($dstRange)

    $ws.Cells.Item($row, 1).Value = $entry.Date
    $ws.Cells.Item($row, 2).Value = $entry.Text
}

# Match the author's final view state: scrolled near the bottom with the
# next empty cell selected.
$ws.Range("B55").Select() | Out-Null
